$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.132.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.901.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4599'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3885'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07875'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9891'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.883.57'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.772'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.044'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07011'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.52%  '

$ws.Range("E17").Value = '  -0.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009902'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9997'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.127.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.323'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.103.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.098'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.901'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '118.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.871'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09320'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8922'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.240'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.321'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.128'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05790'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.165'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02085'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9997'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.659'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5673'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.16%  '

$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.694'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.204'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("E46").Value = '  -1.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07010'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.88%  '

$ws.Range("E48").Value = '  -1.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.551'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("E51").Value = '  -0.07%  '
